$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 12
$ws.Range("H12").Value = 17858050
$ws.Range("I12").Value = 881.55554
$ws.Range("J12").Value = 50000956
$ws.Range("K12").Value = 881.55554
$ws.Range("L12").Value = 50000956
$ws.Range("M12").Value = -711.55554
$ws.Range("N12").Value = -50001296
# Row 17
$ws.Range("H17").Value = 2048.4285
$ws.Range("J17").Value = 2099.8438
$ws.Range("L17").Value = 6299.5314
$ws.Range("N17").Value = -6635.5314
# Row 33
$ws.Range("H33").Value = 11749.92
$ws.Range("I33").Value = 15147.368
$ws.Range("K33").Value = 15147.368
$ws.Range("M33").Value = -14918.368
# Row 38
$ws.Range("H38").Value = 95.14286
$ws.Range("I38").Value = 95.14286
$ws.Range("J38").Value = 0
$ws.Range("K38").Value = 285.42858
$ws.Range("L38").Value = 0
$ws.Range("M38").Value = 86.57141999999999
$ws.Range("N38").ClearContents()
# Row 86
$ws.Range("H86").Value = 9874.333000000001
$ws.Range("J86").Value = 15971
$ws.Range("L86").Value = 15971
$ws.Range("N86").Value = -18217
# Row 89
$ws.Range("H89").Value = 9874.333000000001
$ws.Range("J89").Value = 15971
$ws.Range("L89").Value = 79855
$ws.Range("N89").Value = -91087
# Row 106
$ws.Range("H106").Value = 256028.55
$ws.Range("I106").Value = 373128.88
$ws.Range("J106").Value = 5099.2856
$ws.Range("K106").Value = 373128.88
$ws.Range("L106").Value = 5099.2856
$ws.Range("M106").Value = -372497.88
$ws.Range("N106").Value = -6361.2856
# Row 127
$ws.Range("H127").Value = 1047.7
$ws.Range("I127").Value = 653.8570999999999
$ws.Range("K127").Value = 1961.5713
$ws.Range("M127").Value = 2998.4287
# Row 138
$ws.Range("H138").Value = 15736.232
$ws.Range("J138").Value = 2522.422
$ws.Range("L138").Value = 7567.266
$ws.Range("N138").Value = -17847.266
# Row 141
$ws.Range("H141").Value = 3409.6667
$ws.Range("I141").Value = 3409.6667
$ws.Range("J141").Value = 0
$ws.Range("K141").Value = 10229.0001
$ws.Range("L141").Value = 0
$ws.Range("M141").Value = -5049.000100000001
$ws.Range("N141").ClearContents()

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 45
$ws.Range("H45").Value = 12816.889
$ws.Range("I45").Value = 15620
$ws.Range("J45").Value = 3006
$ws.Range("K45").Value = 15620
$ws.Range("L45").Value = 3006
$ws.Range("M45").Value = -15243
$ws.Range("N45").Value = -3760
# Row 110
$ws.Range("H110").Value = 32720.555
$ws.Range("I110").Value = 32720.555
$ws.Range("K110").Value = 32720.555
$ws.Range("M110").Value = -30675.555
# Row 122
$ws.Range("H122").Value = 2561.9314
$ws.Range("I122").Value = 2458.3674
$ws.Range("K122").Value = 7375.1022
$ws.Range("M122").Value = -4925.1022
# Row 132
$ws.Range("H132").Value = 2208
$ws.Range("I132").Value = 1969.8928
$ws.Range("K132").Value = 5909.678400000001
$ws.Range("M132").Value = -3379.678400000001

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 16
$ws.Range("H16").Value = 0
$ws.Range("J16").Value = 0
$ws.Range("L16").Value = 0
$ws.Range("N16").ClearContents()
# Row 31
$ws.Range("H31").Value = 5000
$ws.Range("J31").Value = 5000
$ws.Range("L31").Value = 5000
$ws.Range("N31").Value = -5504
# Row 86
$ws.Range("H86").Value = 1599.7142
$ws.Range("I86").Value = 1371.2858
$ws.Range("J86").Value = 2285
$ws.Range("K86").Value = 1371.2858
$ws.Range("L86").Value = 2285
$ws.Range("M86").Value = -248.2858000000001
$ws.Range("N86").Value = -4531
# Row 89
$ws.Range("H89").Value = 1599.7142
$ws.Range("I89").Value = 1371.2858
$ws.Range("J89").Value = 2285
$ws.Range("K89").Value = 6856.429
$ws.Range("L89").Value = 11425
$ws.Range("M89").Value = -1240.429
$ws.Range("N89").Value = -22657
# Row 94
$ws.Range("H94").Value = 1123.5
$ws.Range("I94").Value = 1077.6364
$ws.Range("J94").Value = 1179.5555
$ws.Range("K94").Value = 1077.6364
$ws.Range("L94").Value = 1179.5555
$ws.Range("M94").Value = -626.6364000000001
$ws.Range("N94").Value = -2081.5555
# Row 134
$ws.Range("H134").Value = 2325
$ws.Range("I134").Value = 2052.7144
$ws.Range("K134").Value = 6158.1432
$ws.Range("M134").Value = -3623.1432

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 58
$ws.Range("H58").Value = 1361.3704
$ws.Range("I58").Value = 1148.2667
$ws.Range("J58").Value = 1627.75
$ws.Range("K58").Value = 1148.2667
$ws.Range("L58").Value = 1627.75
$ws.Range("M58").Value = -945.2666999999999
$ws.Range("N58").Value = -2033.75
# Row 132
$ws.Range("H132").Value = 1594.44
$ws.Range("I132").Value = 1279.1904
$ws.Range("J132").Value = 3249.5
$ws.Range("K132").Value = 3837.5712
$ws.Range("L132").Value = 9748.5
$ws.Range("M132").Value = -1307.5712
$ws.Range("N132").Value = -14808.5
# Row 134
$ws.Range("H134").Value = 1681.6316
$ws.Range("I134").Value = 1120.037
$ws.Range("J134").Value = 3060.0908
$ws.Range("K134").Value = 3360.111
$ws.Range("L134").Value = 9180.2724
$ws.Range("M134").Value = -825.1109999999999
$ws.Range("N134").Value = -14250.2724
# Row 136
$ws.Range("H136").Value = 1361.3704
$ws.Range("I136").Value = 1148.2667
$ws.Range("J136").Value = 1627.75
$ws.Range("K136").Value = 3444.800099999999
$ws.Range("L136").Value = 4883.25
$ws.Range("M136").Value = -894.8000999999995
$ws.Range("N136").Value = -9983.25

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 131
$ws.Range("H131").Value = 394730.72
$ws.Range("I131").Value = 143714.86
$ws.Range("K131").Value = 431144.58
$ws.Range("M131").Value = -426104.58
# Row 140
$ws.Range("H140").Value = 1279.1818
$ws.Range("I140").Value = 684.8823
$ws.Range("K140").Value = 2054.6469
$ws.Range("M140").Value = 3125.3531
# Row 141
$ws.Range("H141").Value = 3593.9333
$ws.Range("I141").Value = 3593.9333
$ws.Range("K141").Value = 10781.7999
$ws.Range("M141").Value = -5601.7999

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 80
$ws.Range("H80").Value = 76926030
$ws.Range("J80").Value = 3199.8
$ws.Range("L80").Value = 3199.8
$ws.Range("N80").Value = -5195.8
# Row 83
$ws.Range("H83").Value = 76926030
$ws.Range("J83").Value = 3199.8
$ws.Range("L83").Value = 15999
$ws.Range("N83").Value = -25983
# Row 97
$ws.Range("H97").Value = 1534.1666
$ws.Range("I97").Value = 833
$ws.Range("K97").Value = 833
$ws.Range("M97").Value = -337
# Row 134
$ws.Range("H134").Value = 43707.145
$ws.Range("J134").Value = 43707.145
$ws.Range("L134").Value = 131121.435
$ws.Range("N134").Value = -136191.435
# Row 136
$ws.Range("H136").Value = 31525.842
$ws.Range("J136").Value = 31525.842
$ws.Range("L136").Value = 94577.526
$ws.Range("N136").Value = -99677.526

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 25881.467
$ws.Range("I7").Value = 27516.285
$ws.Range("K7").Value = 27516.285
$ws.Range("M7").Value = -27404.285
# Row 22
$ws.Range("H22").Value = 2304.3333
$ws.Range("I22").Value = 638.4
$ws.Range("K22").Value = 638.4
$ws.Range("M22").Value = -343.4
# Row 27
$ws.Range("H27").Value = 2304.3333
$ws.Range("I27").Value = 638.4
$ws.Range("K27").Value = 638.4
$ws.Range("M27").Value = -531.4
# Row 40
$ws.Range("H40").Value = 2060624.6
$ws.Range("I40").Value = 2569.6667
$ws.Range("J40").Value = 6176734.5
$ws.Range("K40").Value = 2569.6667
$ws.Range("L40").Value = 6176734.5
$ws.Range("M40").Value = -2433.6667
$ws.Range("N40").Value = -6177006.5
# Row 88
$ws.Range("H88").Value = 29799.4
$ws.Range("I88").Value = 24749.5
$ws.Range("K88").Value = 24749.5
$ws.Range("M88").Value = -24321.5
# Row 91
$ws.Range("H91").Value = 29799.4
$ws.Range("I91").Value = 24749.5
$ws.Range("K91").Value = 24749.5
$ws.Range("M91").Value = -23267.5
# Row 122
$ws.Range("H122").Value = 6685383.5
$ws.Range("I122").Value = 21017.424
$ws.Range("K122").Value = 63052.272
$ws.Range("M122").Value = -60602.272
# Row 126
$ws.Range("H126").Value = 25881.467
$ws.Range("I126").Value = 27516.285
$ws.Range("K126").Value = 82548.855
$ws.Range("M126").Value = -80078.855
# Row 132
$ws.Range("H132").Value = 2529.3845
$ws.Range("I132").Value = 2154.75
$ws.Range("K132").Value = 6464.25
$ws.Range("M132").Value = -3934.25
# Row 136
$ws.Range("H136").Value = 3201.2654
$ws.Range("J136").Value = 3296.238
$ws.Range("L136").Value = 9888.714
$ws.Range("N136").Value = -14988.714
# Row 138
$ws.Range("H138").Value = 123692
$ws.Range("J138").Value = 123692
$ws.Range("L138").Value = 123692
$ws.Range("N138").Value = -133972
# Row 139
$ws.Range("H139").Value = 88659.836
$ws.Range("I139").Value = 67000
$ws.Range("J139").Value = 99489.75
$ws.Range("K139").Value = 67000
$ws.Range("L139").Value = 99489.75
$ws.Range("M139").Value = -61860
$ws.Range("N139").Value = -109769.75

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 113
$ws.Range("H113").Value = 887.6111
$ws.Range("J113").Value = 999.8570999999999
$ws.Range("L113").Value = 2999.5713
$ws.Range("N113").Value = -7339.5713
# Row 126
$ws.Range("H126").Value = 1719.6471
$ws.Range("I126").Value = 1441.0769
$ws.Range("J126").Value = 2625
$ws.Range("K126").Value = 4323.2307
$ws.Range("M126").Value = -1853.2307
# Row 136
$ws.Range("H136").Value = 1825.25
$ws.Range("I136").Value = 849.1
$ws.Range("J136").Value = 5079.0835
$ws.Range("K136").Value = 2547.3
$ws.Range("L136").Value = 15237.2505
$ws.Range("M136").Value = 2.699999999999818
$ws.Range("N136").Value = -20337.2505
